$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM_PartType-Integrated Stepper")

# Fill in the Manufacturer Part Number column (G) for the two resistors
# that were previously blank. Set G14 first so the new shared strings are
# appended in the same order as the target workbook (CRGCQ0603J470R before
# RMCF0603FT4K70).
$ws.Range("G14").Value = "CRGCQ0603J470R"
$ws.Range("G13").Value = "RMCF0603FT4K70"

# Restore the original cell formatting (writing a value can otherwise
# switch the cell off of the shared "quote prefix" style used by the rest
# of the row).
$ws.Range("F13").Copy()
$ws.Range("G13").PasteSpecial(-4122)
$ws.Range("F14").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection to match the author's final cursor position.
$ws.Range("K11").Select()
